# Update "想去人数" (interest count) figures in the 展览 and 全部类型 sheets
# to match the newly-scraped data (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 sheet (sheet1)
$wsExhibit.Range("F6").Value  = 1264
$wsExhibit.Range("F8").Value  = 7556
$wsExhibit.Range("F10").Value = 109
$wsExhibit.Range("F26").Value = 2599
$wsExhibit.Range("F29").Value = 2795
$wsExhibit.Range("F30").Value = 25
$wsExhibit.Range("F32").Value = 119
$wsExhibit.Range("F34").Value = 634
$wsExhibit.Range("F36").Value = 865
$wsExhibit.Range("F37").Value = 1622
$wsExhibit.Range("F40").Value = 2608

# 全部类型 sheet (sheet4)
$wsAll.Range("F6").Value  = 1264
$wsAll.Range("F7").Value  = 7556
$wsAll.Range("F9").Value  = 109
$wsAll.Range("F26").Value = 2600
$wsAll.Range("F31").Value = 119
$wsAll.Range("F34").Value = 634
$wsAll.Range("F37").Value = 865
$wsAll.Range("F39").Value = 1622
$wsAll.Range("F42").Value = 2608
